$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append new row with the latest mail log entry ---
$logs = $wb.Worksheets.Item("Logs")

$newRow = 54
$logs.Cells.Item($newRow, 1).Value = "Order wijzigen"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Kan ik mijn bestelling nog aanpassen?"
$logs.Cells.Item($newRow, 4).Value = "Bestelling / Levering"
$logs.Cells.Item($newRow, 6).Value = "2025-06-22 22:05:11"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# --- Sheet "Dashboard": category counts changed due to the new log entry ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(9, 1).Value = "Bestelling / Levering"
$dash.Cells.Item(9, 2).Value = 3

$dash.Cells.Item(11, 1).Value = "Overig"
$dash.Cells.Item(11, 2).Value = 2

$dash.Cells.Item(12, 1).Value = "Factuur / Administratie"
$dash.Cells.Item(12, 2).Value = 2

# --- Extend conditional formatting ranges on "Logs" to cover the new row ---
$dFcs = $logs.Range("D2:D54").FormatConditions
for ($i = 1; $i -le $dFcs.Count; $i++) {
    $dFcs.Item($i).ModifyAppliesToRange($logs.Range("D2:D54"))
}

$gFcs = $logs.Range("G2:G54").FormatConditions
for ($i = 1; $i -le $gFcs.Count; $i++) {
    $gFcs.Item($i).ModifyAppliesToRange($logs.Range("G2:G54"))
}
